$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay text so values like "30.300.43" or
# "  -2.04%  " are not coerced into numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.300.43'
$ws.Range("D3").Value = '1.878.82'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '237.11'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '0.4810'
$ws.Range("E7").Value = '  -2.39%  '
$ws.Range("D8").Value = '0.2884'
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("D9").Value = '0.06585'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("D10").Value = '1.879.51'
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").Value = '16.93'
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("D12").Value = '0.07383'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '5.196'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").Value = '87.86'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = '0.6593'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '30.264.16'
$ws.Range("D17").Value = '13.53'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").Value = '0.000007715'
$ws.Range("E19").Value = '  -2.68%  '
$ws.Range("D20").Value = '5.455'
$ws.Range("E20").Value = '  +1.74%  '
$ws.Range("D21").Value = '2.140.93'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '194.15'
$ws.Range("E23").Value = '  -4.73%  '
$ws.Range("D24").Value = '6.178'
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("D25").Value = '9.433'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").Value = '165.16'
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("D27").Value = '18.24'
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").Value = '1.926'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").Value = '1.445'
$ws.Range("E29").Value = '  -2.52%  '
$ws.Range("D30").Value = '4.271'
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("D31").Value = '0.09140'
$ws.Range("D32").Value = '4.045'
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").Value = '0.05055'
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("D34").Value = '0.7381'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = '1.136'
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("D36").Value = '2.711'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '0.01849'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '2.631'
$ws.Range("E38").Value = '  -3.23%  '
$ws.Range("D39").Value = '0.9151'
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").Value = '2.074'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").Value = '106.34'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.887'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4321'
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").Value = '7.661'
$ws.Range("E45").Value = '  +0.35%  '
$ws.Range("D46").Value = '0.1347'
$ws.Range("E46").Value = '  -3.06%  '
$ws.Range("D47").Value = '1.586'
$ws.Range("E47").Value = '  +10.29%  '
$ws.Range("D48").Value = '65.36'
$ws.Range("E48").Value = '  -10.18%  '
$ws.Range("D49").Value = '8.877'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").Value = '34.14'
$ws.Range("E50").Value = '  -3.31%  '
$ws.Range("E51").Value = '  -2.63%  '
